$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 105, pushing the existing rows 105-118 down to 106-119
$ws.Rows.Item(105).Insert()

# Populate the newly inserted row 105 with the new weekly price record
$ws.Range("A105").Value = 4
$ws.Range("B105").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C105").Value = "Los Lagos"
$ws.Range("D105").Value = 44446
$ws.Range("E105").Value = 10
$ws.Range("F105").Value = 100112032
$ws.Range("G105").Value = "Zapallo italiano"
$ws.Range("H105").Value = "Sin especificar"
$ws.Range("I105").Value = "Primera"
$ws.Range("J105").Value = 200
$ws.Range("K105").Value = 20000
$ws.Range("L105").Value = 21000
$ws.Range("M105").Value = 20500
$ws.Range("N105").Value = "$/caja 50 unidades"
$ws.Range("O105").Value = "Región de Arica y Parinacota"
$ws.Range("P105").Value = 410
$ws.Range("Q105").Value = 50
$ws.Range("R105").Value = "Hortaliza"
